# Rename the inline picture shapes embedded in the document's headers and
# footers. The Pearson Edexcel logo (currently reported as "image1.png")
# becomes "image2.png", and the BTEC logo (currently reported as
# "image2.jpg") becomes "image1.jpg" - the two images effectively swap
# their shape-name numbering.
#
# The pictures live in the header/footer stories, not the main body, so we
# reach them through Sections(1).Headers(...)/Footers(...) rather than
# $d.InlineShapes (which only covers the main document story).

$d = $word.ActiveDocument
$sec = $d.Sections(1)

# First-page header holds the BTEC logo.
$btecHeader = $sec.Headers(2)
if ($btecHeader.Exists -and $btecHeader.Range.InlineShapes.Count -ge 1) {
    $btecShape = $btecHeader.Range.InlineShapes.Item(1)
    $btecShape.Name = "image1.jpg"
}

# Default (primary) footer holds a Pearson logo.
$pearsonFooter1 = $sec.Footers(1)
if ($pearsonFooter1.Exists -and $pearsonFooter1.Range.InlineShapes.Count -ge 1) {
    $pearsonShape1 = $pearsonFooter1.Range.InlineShapes.Item(1)
    $pearsonShape1.Name = "image2.png"
}

# First-page footer holds a Pearson logo too.
$pearsonFooter2 = $sec.Footers(2)
if ($pearsonFooter2.Exists -and $pearsonFooter2.Range.InlineShapes.Count -ge 1) {
    $pearsonShape2 = $pearsonFooter2.Range.InlineShapes.Item(1)
    $pearsonShape2.Name = "image2.png"
}
